$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet tab from "Java" to "AWS"
$ws.Name = "AWS"

# 2. Update row 2 (Starter Lambda description / reformatted Lambda code)
#    B2 is written first so the rebuilt shared-string table keeps the code
#    string ahead of the description string (matches target string order)
$ws.Range("B2").Value = "exports.handler = async (event, context, callback) => {`n    const hasError = event['queryStringParameters']['myErrorParam'];`n    if (hasError === 'yes') {`n        callback(new Error('My error message'));`n    }`n    else {`n        const response = {`n            statusCode: 200,`n            headers: {`"Access-Control-Allow-Origin`": `"*`"},`n            body: JSON.stringify({success: true}),`n            isBase64Encoded: false`n        };`n        callback(null, response);`n}"
$ws.Range("A2").Value = "Starter Lambda function compatible with API Gateway v1.0"

# 3. Update row 3 (Batch script description / batch file content) - unchanged content,
#    rewritten here (B3 before A3) so the shared-string table ends up in the
#    same order as the target file
$ws.Range("B3").Value = ":: This batch file redeploys an existing lambda function`n:: Usage: deploy-existing.bat getAllProducts`necho off`necho WARNING: this will delete any index.js or index.zip you have in the current directory!`npause`nset /p toDeploy=Enter lambda name (without the .js): `npowershell -Command `"(gc %toDeploy%.js) -replace './helpers', '/opt/lambdas/helpers' | Out-File -encoding ASCII index.js`"`npowershell `"Compress-Archive index.js index.zip`"`naws lambda update-function-code --function-name %toDeploy% --zip-file fileb://index.zip`ndel index.js`ndel index.zip"
$ws.Range("A3").Value = "Batch script to update an existing Lambda function that uses a Lambda Layer. Assumes local layer name is ./helpers and layer on AWS Lambda is located at /opt/lambdas/helpers"

# 4. Update the sheet view so the top-left cell and selection are A2 instead of A3
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Done"
